# Update column G ("K" = strikeouts) values for rows 2-36 in the
# kuhl_chad.xlsx "save_data" sheet. These values are recalculated/regenerated
# from box-score strikeout totals (K) rather than the previous "Strike#"
# derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 3
    8  = 1
    9  = 1
    10 = 2
    11 = 1
    12 = 4
    13 = 2
    14 = 2
    15 = 2
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 2
    21 = 2
    22 = 1
    23 = 4
    24 = 0
    25 = 2
    26 = 4
    27 = 1
    28 = 1
    29 = 3
    30 = 7
    31 = 1
    32 = 3
    33 = 2
    34 = 0
    35 = 2
    36 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
